# Bug: Illegal sheet symbols; incomplete cell squares
# Close #11.
#
# - Add a new worksheet "Sheet2" after "Sheet1" containing a formula that
#   sums Sheet1's A2:A11 range (referencing the sheet by name, since that
#   was previously illegal / produced "incomplete" results when computed
#   locally on Sheet1 only).
# - Make Sheet2 the active sheet/tab, with C4 selected.
# - Clear the "tabSelected" / previous selection state on Sheet1 and
#   select E12 there instead.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# New formula on Sheet2 summing Sheet1's running totals.
$ws2.Range("B2").Formula = "=SUM(Sheet1!A2:A11)"

# Update Sheet1's selection/view state (it is no longer the active tab).
$ws1.Range("E12").Select()

# Make Sheet2 the active sheet and select C4 there.
$ws2.Activate()
$ws2.Range("C4").Select()
